# Apply updated crypto price/volume figures per the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $text) {
    # Force the cell to keep/become a plain text value, even when the
    # string looks like a number (e.g. "1.000"), then restore the
    # cell formatting so no visible style change is introduced.
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.NumberFormat = "General"
    $range.Style = "Normal"
}

$ws.Range("D2").Value = "27.715.00"
$ws.Range("E2").Value = "  +0.22%  "

$ws.Range("D3").Value = "1.891.54"
$ws.Range("E3").Value = "  +1.19%  "

Set-TextValue $ws.Range("D4") "1.000"

Set-TextValue $ws.Range("D5") "313.44"
$ws.Range("E5").Value = "  +0.16%  "

Set-TextValue $ws.Range("D6") "1.000"
$ws.Range("E6").Value = "  -1.21%  "

Set-TextValue $ws.Range("D7") "0.4851"
$ws.Range("E7").Value = "  +0.45%  "

Set-TextValue $ws.Range("D8") "0.3795"
$ws.Range("E8").Value = "  -0.45%  "

Set-TextValue $ws.Range("D9") "0.07336"
$ws.Range("E9").Value = "  -0.46%  "

Set-TextValue $ws.Range("D10") "0.9187"
$ws.Range("E10").Value = "  -1.84%  "

Set-TextValue $ws.Range("D11") "20.53"
$ws.Range("E11").Value = "  -1.71%  "

$ws.Range("E12").Value = "  -1.53%  "

$ws.Range("D13").Value = "1.882.93"
$ws.Range("E13").Value = "  +0.42%  "

Set-TextValue $ws.Range("D14") "5.463"
$ws.Range("E14").Value = "  -0.10%  "

Set-TextValue $ws.Range("D15") "6.598"
$ws.Range("E15").Value = "  +0.27%  "

Set-TextValue $ws.Range("D16") "90.99"
$ws.Range("E16").Value = "  +0.62%  "

Set-TextValue $ws.Range("D17") "1.001"
$ws.Range("E17").Value = "  -1.29%  "

Set-TextValue $ws.Range("D18") "0.000008801"
$ws.Range("E18").Value = "  -0.63%  "

$ws.Range("E19").Value = "  -1.05%  "

$ws.Range("D20").Value = "27.737.83"
$ws.Range("E20").Value = "  +0.06%  "

$ws.Range("E21").Value = "  -1.54%  "

Set-TextValue $ws.Range("D22") "5.117"
$ws.Range("E22").Value = "  +0.00%  "

$ws.Range("D23").Value = "2.127.25"
$ws.Range("E23").Value = "  +0.88%  "

Set-TextValue $ws.Range("D24") "10.81"

Set-TextValue $ws.Range("D25") "1.910"
$ws.Range("E25").Value = "  -1.82%  "

Set-TextValue $ws.Range("D26") "153.27"
$ws.Range("E26").Value = "  -2.09%  "

Set-TextValue $ws.Range("D27") "18.37"
$ws.Range("E27").Value = "  -0.94%  "

Set-TextValue $ws.Range("D28") "2.130"
$ws.Range("E28").Value = "  +3.46%  "

Set-TextValue $ws.Range("D29") "115.84"
$ws.Range("E29").Value = "  +0.16%  "

Set-TextValue $ws.Range("D30") "4.902"
$ws.Range("E30").Value = "  -1.14%  "

Set-TextValue $ws.Range("D31") "0.08935"
$ws.Range("E31").Value = "  +0.12%  "

Set-TextValue $ws.Range("D32") "3.150"
$ws.Range("E32").Value = "  -5.64%  "

Set-TextValue $ws.Range("D33") "1.222"
$ws.Range("E33").Value = "  +0.60%  "

Set-TextValue $ws.Range("D34") "0.7631"
$ws.Range("E34").Value = "  +0.24%  "

Set-TextValue $ws.Range("D35") "4.632"
$ws.Range("E35").Value = "  +0.29%  "

Set-TextValue $ws.Range("D36") "0.02035"
$ws.Range("E36").Value = "  -0.23%  "

Set-TextValue $ws.Range("D37") "2.543"
$ws.Range("E37").Value = "  -6.24%  "

Set-TextValue $ws.Range("D38") "1.092"
$ws.Range("E38").Value = "  -3.98%  "

Set-TextValue $ws.Range("D39") "0.05265"
$ws.Range("E39").Value = "  -2.25%  "

Set-TextValue $ws.Range("D40") "0.5461"
$ws.Range("E40").Value = "  -3.25%  "

Set-TextValue $ws.Range("D41") "2.981"
$ws.Range("E41").Value = "  -0.14%  "

Set-TextValue $ws.Range("D42") "6.940"
$ws.Range("E42").Value = "  -1.68%  "

Set-TextValue $ws.Range("D43") "0.1520"
$ws.Range("E43").Value = "  -0.66%  "

Set-TextValue $ws.Range("D44") "8.333"
$ws.Range("E44").Value = "  -2.66%  "

Set-TextValue $ws.Range("D45") "109.76"
$ws.Range("E45").Value = "  +4.59%  "

Set-TextValue $ws.Range("D46") "10.60"
$ws.Range("E46").Value = "  -1.64%  "

Set-TextValue $ws.Range("D47") "0.4785"
$ws.Range("E47").Value = "  -2.50%  "

Set-TextValue $ws.Range("D48") "1.000"
$ws.Range("E48").Value = "  -1.31%  "

Set-TextValue $ws.Range("D49") "1.635"
$ws.Range("E49").Value = "  -2.26%  "

Set-TextValue $ws.Range("D50") "67.29"
$ws.Range("E50").Value = "  -0.72%  "

Set-TextValue $ws.Range("D51") "0.06048"
$ws.Range("E51").Value = "  -0.89%  "
